# Refresh the crypto price/volume snapshot (Price column D, Volume(1h)
# column E) to the latest scraped values. Some Price strings (e.g.
# "204.10", "1.50", "0.990") are plain decimal-looking text, so the cell's
# NumberFormat is forced to Text ("@") right before the assignment --
# otherwise Excel's COM layer would silently coerce them to real numbers
# and drop the trailing zero / literal formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.779.30"
$ws.Range("E2").Value = "  -1.70%  "
$ws.Range("D3").Value = "1.550.90"
$ws.Range("E3").Value = "  -1.57%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "204.10"
$ws.Range("E5").Value = "  -1.74%  "
$ws.Range("E6").Value = "  -1.83%  "
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("E8").Value = "  -4.27%  "
$ws.Range("E9").Value = "  -1.19%  "
$ws.Range("E10").Value = "  -1.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0858"
$ws.Range("E11").Value = "  -0.88%  "
$ws.Range("D12").Value = "1.773.07"
$ws.Range("E12").Value = "  -1.46%  "
$ws.Range("D13").Value = "1.559.28"
$ws.Range("E13").Value = "  -0.75%  "
$ws.Range("E14").Value = "  -2.80%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.509"
$ws.Range("E15").Value = "  -2.19%  "
$ws.Range("D16").Value = "26.768.91"
$ws.Range("E16").Value = "  -1.69%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.08"
$ws.Range("E17").Value = "  -2.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "214.09"
$ws.Range("E18").Value = "  -0.92%  "
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("D20").Value = "0.0₃0679"
$ws.Range("E20").Value = "  -1.09%  "
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.07"
$ws.Range("E22").Value = "  -1.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.06"
$ws.Range("E23").Value = "  -3.91%  "
$ws.Range("E24").Value = "  -0.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.68"
$ws.Range("E25").Value = "  -0.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.51"
$ws.Range("E26").Value = "  -2.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.85"
$ws.Range("E27").Value = "  -0.94%  "
$ws.Range("E28").Value = "  +0.26%  "
$ws.Range("E29").Value = "  -2.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0462"
$ws.Range("E30").Value = "  -0.54%  "
$ws.Range("E31").Value = "  -3.10%  "
$ws.Range("E32").Value = "  -0.58%  "
$ws.Range("D33").Value = "1.364.86"
$ws.Range("E33").Value = "  -2.57%  "
$ws.Range("E34").Value = "  -1.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.50"
$ws.Range("E35").Value = "  -4.60%  "
$ws.Range("E36").Value = "  -0.39%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.920"
$ws.Range("E37").Value = "  -2.76%  "
$ws.Range("E38").Value = "  -2.36%  "
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.801"
$ws.Range("E40").Value = "  -2.27%  "
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.990"
$ws.Range("E42").Value = "  -0.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.55"
$ws.Range("E43").Value = "  +3.41%  "
$ws.Range("E44").Value = "  +0.19%  "
$ws.Range("E45").Value = "  -3.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "62.87"
$ws.Range("E46").Value = "  -1.66%  "
$ws.Range("E47").Value = "  -2.14%  "
$ws.Range("D48").Value = "1.686.76"
$ws.Range("E48").Value = "  -1.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "86.03"
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("E50").Value = "  +3.49%  "
$ws.Range("D51").Value = "0.0₇0977"
$ws.Range("E51").Value = "  +0.36%  "
